$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update UsdRate column (D) from 8.15 to 9.25 for all data rows
$ws.Range("D2:D10").Value = 9.25

# Add new HomePage column (E): header then the Luxeon row's URL
$ws.Range("E1").Value = "HomePage"
$ws.Range("E2").Value = "http://luxeon.ua/"

# Add new Description column (F): header then the Luxeon row's description
$ws.Range("F1").Value = "Description"
$ws.Range("F2").Value = "Під торговою маркою LUXEON з 2000 року виробляється широкий спектр обладнання - побутові та промислові системи стабілізації напруги і безперебійного живлення, комплектуючі для персональних комп'ютерів, акустичні системи. З весни 2008 року на прилавках магазинів з'явилася побутова техніка, вироблена під ТМ LUXEON - телевізори, ДВД- програвачі, СВЧ - печі, холодильники і кондиціонери.<br/>Роблячи вибір на користь товарів під ТМ «LUXEON» Ви можете бути впевнені, що отримаєте якісну і надійну техніку, яка завжди буде відповідати Вашим очікуванням!"

# Set explicit column widths to match target layout
$ws.Columns.Item(5).ColumnWidth = 16.5
$ws.Columns.Item(6).ColumnWidth = 13.333333333333334

# Update the selection: previously E3 was selected, then B2 was ctrl-clicked (making it active)
$ws.Range("E3").Select()
$ws.Range("B2").Activate()
